$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ C = 0.5453776865001148;   E = -0.7976031984000098 }
    3  = @{ C = -6.170514117037273;   E = -8.396348489509153 }
    4  = @{ C = -0.2290082001396909;  E = -4.327930935900004 }
    5  = @{ C = 1.406827509327035;    E = 2.015050062499957 }
    6  = @{ C = 1.153683074671208;    E = 3.648892256099945 }
    7  = @{ C = 0.2186142574756467;   E = 0.4006004000999708 }
    8  = @{ C = -0.8522658067264599;  E = -3.551690943899999 }
    9  = @{ C = -0.2262139320475365;  E = -0.7976031983999876 }
    10 = @{ C = -0.6258176826215101;  E = -0.3994003999000073 }
    11 = @{ C = 0.3239252862367037;   E = 1.609625625600009 }
    12 = @{ C = 0.5738128002843901;   E = -0.3994003999000184 }
    13 = @{ C = -0.4781004700720293;  E = 0.8024032015999882 }
    14 = @{ C = -1.197849743493773;   E = -3.161804390399992 }
    15 = @{ C = 1.064698711638945;    E = -2.540956581357878 }
    16 = @{ C = -1.524103236349472;   E = -1.240907591477092 }
    17 = @{ C = 0.2854413827033664;   E = -0.2470349027347551 }
    18 = @{ C = 1.118108578853261;    E = 1.532721825047534 }
    19 = @{ C = -1.490505436658163;   E = -0.3349088112516219 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
